$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 155; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $v = $cell.Value2
    if ($v -is [double]) {
        $cell.Value2 = -$v
    }
}
